# Ajuste na pontuação das tabelas
# Divide each value in column Q (rows 2-163) by 24 and round to 2 decimals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 163; $r++) {
    $cell = $ws.Cells.Item($r, 17)  # Column Q = 17
    $oldValue = $cell.Value2
    if ($oldValue -ne $null) {
        $newValue = [Math]::Round([double]$oldValue / 24, 2)
        $cell.Value = $newValue
    }
}
